# The deck currently has the "Integral" theme applied to the slide master
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) while the original default
# "Office Theme" palette only survives, unused, as the notes master's theme.
# The authored edit swaps the two: the presentation's live theme becomes the
# plain "Office Theme" palette again. Reproduce that by rewriting each of
# the twelve standard theme colors on the slide master's color scheme to the
# "Office Theme" values (PowerPoint packs RGB() as 0x00BBGGRR).

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$officeThemeColors = @(
    0x000000,  # 1  dk1       000000
    0xFFFFFF,  # 2  lt1       FFFFFF
    0x6A5444,  # 3  dk2       44546A
    0xE6E6E7,  # 4  lt2       E7E6E6
    0xD59B5B,  # 5  accent1   5B9BD5
    0x317DED,  # 6  accent2   ED7D31
    0xA5A5A5,  # 7  accent3   A5A5A5
    0x00C0FF,  # 8  accent4   FFC000
    0xC47244,  # 9  accent5   4472C4
    0x47AD70,  # 10 accent6   70AD47
    0xC16305,  # 11 hlink     0563C1
    0x724F95   # 12 folHlink  954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $cs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
